$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data cells for row 2
$ws.Range("B2").Value = $false
$ws.Range("E2").Value = 77.349997999999999
$ws.Range("F2").Value = -1.0869590792838919
$ws.Range("G2").Value = $false

# New row 3
$ws.Range("C3").Value = 9891.2999999999993

# Column width adjustments (closest representable values given the
# engine's internal pixel-rounding of character widths)
$ws.Columns.Item(5).ColumnWidth = 8.92
$ws.Columns.Item(6).ColumnWidth = 11.67
